{"js": "// 1) Trim the trailing \", \u0430 \u043f\u0440\u043e\u0441\u0442\u0440\u0430\u043d\u0441\u0442\u0432\u0435\u043d\u043d\u0430\u044f \u0441\u043b\u043e\u0436\u043d\u043e\u0441\u0442\u044c O(n)\" clause from the\n//    recursion-complexity sentence, leaving the sentence ending at \"O(2^n).\"\nconst oldTail = \", \u0430 \u043f\u0440\u043e\u0441\u0442\u0440\u0430\u043d\u0441\u0442\u0432\u0435\u043d\u043d\u0430\u044f \u0441\u043b\u043e\u0436\u043d\u043e\u0441\u0442\u044c O(n).\";\nconst searchResults = context.document.body.search(oldTail, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\".\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Remove the trailing empty paragraph and the \"\u041c\u043e\u0436\u043d\u043e \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u0432\u044b\u0432\u043e\u0434...\"\n//    paragraph that followed the iterative-approach paragraph at the very\n//    end of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet conclusionIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"\u041c\u043e\u0436\u043d\u043e \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u0432\u044b\u0432\u043e\u0434 \u043e \u0442\u043e\u043c, \u0447\u0442\u043e \u0440\u0435\u043a\u0443\u0440\u0441\u0438\u0432\u043d\u044b\u0439 \u043f\u043e\u0434\u0445\u043e\u0434\") !== -1) {\n    conclusionIndex = i;\n    break;\n  }\n}\n\nif (conclusionIndex !== -1) {\n  items[conclusionIndex].delete();\n  if (conclusionIndex - 1 >= 0 && items[conclusionIndex - 1].text === \"\") {\n    items[conclusionIndex - 1].delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Trim the trailing \", \u0430 \u043f\u0440\u043e\u0441\u0442\u0440\u0430\u043d\u0441\u0442\u0432\u0435\u043d\u043d\u0430\u044f \u0441\u043b\u043e\u0436\u043d\u043e\u0441\u0442\u044c O(n)\" clause from the\n#    recursion-complexity sentence, leaving the sentence ending at \"O(2^n).\"\n$find = $d.Content.Find\n$find.Text = \", \u0430 \u043f\u0440\u043e\u0441\u0442\u0440\u0430\u043d\u0441\u0442\u0432\u0435\u043d\u043d\u0430\u044f \u0441\u043b\u043e\u0436\u043d\u043e\u0441\u0442\u044c O(n).\"\n$find.Replacement.Text = \".\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n# 2) Remove the trailing empty paragraph and the \"\u041c\u043e\u0436\u043d\u043e \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u0432\u044b\u0432\u043e\u0434...\"\n#    paragraph that followed the iterative-approach paragraph at the very\n#    end of the document.\n$count = $d.Paragraphs.Count\n$conclusionIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text -like \"*\u041c\u043e\u0436\u043d\u043e \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u0432\u044b\u0432\u043e\u0434 \u043e \u0442\u043e\u043c, \u0447\u0442\u043e \u0440\u0435\u043a\u0443\u0440\u0441\u0438\u0432\u043d\u044b\u0439 \u043f\u043e\u0434\u0445\u043e\u0434*\") {\n    $conclusionIndex = $i\n    break\n  }\n}\n\nif ($conclusionIndex -gt 0) {\n  $d.Paragraphs.Item($conclusionIndex).Range.Delete()\n  if ($conclusionIndex - 1 -ge 1) {\n    $prev = $d.Paragraphs.Item($conclusionIndex - 1)\n    if ($prev.Range.Text.Trim().Length -eq 0) {\n      $prev.Range.Delete()\n    }\n  }\n}\n"}
